$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Fix the search test case results from "N" to "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Update the active selection to reflect the fixed cell
$ws.Range("C3").Select()
